$wb = $excel.ActiveWorkbook

# --- Sheet: ALC (index 1) ---
$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(17, 8).Value = 2624.0857
$ws.Cells.Item(17, 10).Value = 2264.1516
$ws.Cells.Item(17, 12).Value = 6792.4548
$ws.Cells.Item(17, 14).Value = -7128.4548
$ws.Cells.Item(32, 8).Value = 2272.6365
$ws.Cells.Item(32, 9).Value = 2733
$ws.Cells.Item(32, 10).Value = 2100
$ws.Cells.Item(32, 11).Value = 2733
$ws.Cells.Item(32, 12).Value = 2100
$ws.Cells.Item(32, 13).Value = -2407
$ws.Cells.Item(32, 14).Value = -2752
$ws.Cells.Item(51, 8).Value = 6200.2
$ws.Cells.Item(51, 9).Value = 5500.5
$ws.Cells.Item(51, 10).Value = 6666.6665
$ws.Cells.Item(51, 11).Value = 5500.5
$ws.Cells.Item(51, 12).Value = 6666.6665
$ws.Cells.Item(51, 13).Value = -5016.5
$ws.Cells.Item(51, 14).Value = -7634.6665
$ws.Cells.Item(100, 8).Value = 1499.1
$ws.Cells.Item(100, 9).Value = 743.2
$ws.Cells.Item(100, 11).Value = 743.2
$ws.Cells.Item(100, 13).Value = -202.2
$ws.Cells.Item(129, 8).Value = 1113.8
$ws.Cells.Item(129, 10).Value = 1157.4445
$ws.Cells.Item(129, 12).Value = 3472.3335
$ws.Cells.Item(129, 14).Value = -13472.3335
$ws.Cells.Item(132, 8).Value = 1479.6757
$ws.Cells.Item(132, 9).Value = 1295.5758
$ws.Cells.Item(132, 11).Value = 3886.7274
$ws.Cells.Item(132, 13).Value = -1356.7274
$ws.Cells.Item(137, 8).Value = 2409.4167
$ws.Cells.Item(137, 9).Value = 1807.2858
$ws.Cells.Item(137, 10).Value = 2657.353
$ws.Cells.Item(137, 11).Value = 5421.857400000001
$ws.Cells.Item(137, 12).Value = 7972.059
$ws.Cells.Item(137, 13).Value = -2871.857400000001
$ws.Cells.Item(137, 14).Value = -13072.059
$ws.Cells.Item(138, 8).Value = 3070.5151
$ws.Cells.Item(138, 9).Value = 3578.0625
$ws.Cells.Item(138, 10).Value = 2592.8235
$ws.Cells.Item(138, 11).Value = 10734.1875
$ws.Cells.Item(138, 12).Value = 7778.470499999999
$ws.Cells.Item(138, 13).Value = -5594.1875
$ws.Cells.Item(138, 14).Value = -18058.4705

# --- Sheet: ARM (index 2) ---
$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(22, 8).Value = 3499.5
$ws.Cells.Item(22, 9).Value = 3499.5
$ws.Cells.Item(22, 11).Value = 3499.5
$ws.Cells.Item(22, 13).Value = -3200.5
$ws.Cells.Item(28, 8).Value = 21666.334
$ws.Cells.Item(28, 9).Value = 21666.334
$ws.Cells.Item(28, 11).Value = 21666.334
$ws.Cells.Item(28, 13).Value = -21474.334
$ws.Cells.Item(32, 8).Value = 4538.34
$ws.Cells.Item(32, 9).Value = 3643.682
$ws.Cells.Item(32, 10).Value = 11099.167
$ws.Cells.Item(32, 11).Value = 3643.682
$ws.Cells.Item(32, 12).Value = 11099.167
$ws.Cells.Item(32, 13).Value = -3356.682
$ws.Cells.Item(32, 14).Value = -11673.167
$ws.Cells.Item(41, 8).Value = 30531
$ws.Cells.Item(41, 9).Value = 4000
$ws.Cells.Item(41, 10).Value = 57062
$ws.Cells.Item(41, 11).Value = 4000
$ws.Cells.Item(41, 12).Value = 57062
$ws.Cells.Item(41, 13).Value = -3586
$ws.Cells.Item(41, 14).Value = -57890
$ws.Cells.Item(45, 8).Value = 1693.8572
$ws.Cells.Item(45, 9).Value = 941
$ws.Cells.Item(45, 11).Value = 941
$ws.Cells.Item(45, 13).Value = -564
$ws.Cells.Item(97, 8).Value = 1066.1111
$ws.Cells.Item(97, 9).Value = 1037.3529
$ws.Cells.Item(97, 11).Value = 1037.3529
$ws.Cells.Item(97, 13).Value = -541.3529000000001
$ws.Cells.Item(99, 8).Value = 21666.334
$ws.Cells.Item(99, 9).Value = 21666.334
$ws.Cells.Item(99, 11).Value = 21666.334
$ws.Cells.Item(99, 13).Value = -18671.334
$ws.Cells.Item(110, 8).Value = 236.18182
$ws.Cells.Item(110, 9).Value = 216.11111
$ws.Cells.Item(110, 11).Value = 216.11111
$ws.Cells.Item(110, 13).Value = 1828.88889
$ws.Cells.Item(122, 8).Value = 2267.1428
$ws.Cells.Item(122, 9).Value = 1850.5555
$ws.Cells.Item(122, 10).Value = 4766.6665
$ws.Cells.Item(122, 11).Value = 5551.666499999999
$ws.Cells.Item(122, 12).Value = 14299.9995
$ws.Cells.Item(122, 13).Value = -3101.666499999999
$ws.Cells.Item(122, 14).Value = -19199.9995
$ws.Cells.Item(132, 8).Value = 1306.8572
$ws.Cells.Item(132, 9).Value = 1048.7333
$ws.Cells.Item(132, 10).Value = 2855.6
$ws.Cells.Item(132, 11).Value = 3146.199900000001
$ws.Cells.Item(132, 12).Value = 8566.799999999999
$ws.Cells.Item(132, 13).Value = -616.1999000000005
$ws.Cells.Item(132, 14).Value = -13626.8

# --- Sheet: BSM (index 3) ---
$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(94, 8).Value = 745.1429000000001
$ws.Cells.Item(94, 9).Value = 814.3333
$ws.Cells.Item(94, 10).Value = 330
$ws.Cells.Item(94, 11).Value = 814.3333
$ws.Cells.Item(94, 12).Value = 330
$ws.Cells.Item(94, 13).Value = -363.3333
$ws.Cells.Item(94, 14).Value = -1232

# --- Sheet: CRP (index 4) ---
$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(31, 8).Value = 2306.2856
$ws.Cells.Item(31, 9).Value = 1974.3846
$ws.Cells.Item(31, 11).Value = 1974.3846
$ws.Cells.Item(31, 13).Value = -1679.3846
$ws.Cells.Item(34, 8).Value = 2306.2856
$ws.Cells.Item(34, 9).Value = 1974.3846
$ws.Cells.Item(34, 11).Value = 1974.3846
$ws.Cells.Item(34, 13).Value = -1772.3846
$ws.Cells.Item(58, 8).Value = 1893669.2
$ws.Cells.Item(58, 9).Value = 3108043.8
$ws.Cells.Item(58, 11).Value = 3108043.8
$ws.Cells.Item(58, 13).Value = -3107840.8
$ws.Cells.Item(59, 8).Value = 18040
$ws.Cells.Item(59, 10).Value = 18040
$ws.Cells.Item(59, 12).Value = 18040
$ws.Cells.Item(59, 14).Value = -20330
$ws.Cells.Item(86, 8).Value = 3660.75
$ws.Cells.Item(86, 9).Value = 3252.8
$ws.Cells.Item(86, 10).Value = 4340.6665
$ws.Cells.Item(86, 11).Value = 3252.8
$ws.Cells.Item(86, 12).Value = 4340.6665
$ws.Cells.Item(86, 13).Value = -2129.8
$ws.Cells.Item(86, 14).Value = -6586.6665
$ws.Cells.Item(89, 8).Value = 3660.75
$ws.Cells.Item(89, 9).Value = 3252.8
$ws.Cells.Item(89, 10).Value = 4340.6665
$ws.Cells.Item(89, 11).Value = 16264
$ws.Cells.Item(89, 12).Value = 21703.3325
$ws.Cells.Item(89, 13).Value = -10648
$ws.Cells.Item(89, 14).Value = -32935.3325
$ws.Cells.Item(99, 8).Value = 2452.9167
$ws.Cells.Item(99, 9).Value = 1491.875
$ws.Cells.Item(99, 11).Value = 1491.875
$ws.Cells.Item(99, 13).Value = 6.125
$ws.Cells.Item(122, 8).Value = 3778.8
$ws.Cells.Item(122, 9).Value = 2928.1
$ws.Cells.Item(122, 10).Value = 5480.2
$ws.Cells.Item(122, 11).Value = 8784.299999999999
$ws.Cells.Item(122, 12).Value = 16440.6
$ws.Cells.Item(122, 13).Value = -6334.299999999999
$ws.Cells.Item(122, 14).Value = -21340.6
$ws.Cells.Item(126, 8).Value = 2452.9167
$ws.Cells.Item(126, 9).Value = 1491.875
$ws.Cells.Item(126, 11).Value = 4475.625
$ws.Cells.Item(126, 13).Value = -2005.625
$ws.Cells.Item(132, 8).Value = 2340.6487
$ws.Cells.Item(132, 9).Value = 1705.3572
$ws.Cells.Item(132, 11).Value = 5116.071599999999
$ws.Cells.Item(132, 13).Value = -2586.071599999999
$ws.Cells.Item(136, 8).Value = 1893669.2
$ws.Cells.Item(136, 9).Value = 3108043.8
$ws.Cells.Item(136, 11).Value = 9324131.399999999
$ws.Cells.Item(136, 13).Value = -9321581.399999999

# --- Sheet: GSM (index 6) ---
$ws = $wb.Worksheets.Item(6)
$ws.Cells.Item(97, 8).Value = 1991.7142
$ws.Cells.Item(97, 9).Value = 1898.7273
$ws.Cells.Item(97, 10).Value = 2332.6667
$ws.Cells.Item(97, 11).Value = 1898.7273
$ws.Cells.Item(97, 12).Value = 2332.6667
$ws.Cells.Item(97, 13).Value = -1402.7273
$ws.Cells.Item(97, 14).Value = -3324.6667
$ws.Cells.Item(102, 8).Value = 3905.9614
$ws.Cells.Item(102, 9).Value = 4779.467
$ws.Cells.Item(102, 11).Value = 4779.467
$ws.Cells.Item(102, 13).Value = -3157.467
$ws.Cells.Item(126, 8).Value = 1259198.9
$ws.Cells.Item(126, 9).Value = 3270510.8
$ws.Cells.Item(126, 10).Value = 38045.355
$ws.Cells.Item(126, 11).Value = 9811532.399999999
$ws.Cells.Item(126, 12).Value = 114136.065
$ws.Cells.Item(126, 13).Value = -9809062.399999999
$ws.Cells.Item(126, 14).Value = -119076.065

# --- Sheet: LTW (index 7) ---
$ws = $wb.Worksheets.Item(7)
$ws.Cells.Item(7, 8).Value = 4074.9333
$ws.Cells.Item(7, 9).Value = 3306.875
$ws.Cells.Item(7, 10).Value = 4952.7144
$ws.Cells.Item(7, 11).Value = 3306.875
$ws.Cells.Item(7, 12).Value = 4952.7144
$ws.Cells.Item(7, 13).Value = -3194.875
$ws.Cells.Item(7, 14).Value = -5176.7144
$ws.Cells.Item(22, 8).Value = 1943.2
$ws.Cells.Item(22, 10).Value = 2020.2222
$ws.Cells.Item(22, 12).Value = 2020.2222
$ws.Cells.Item(22, 14).Value = -2610.2222
$ws.Cells.Item(27, 8).Value = 1943.2
$ws.Cells.Item(27, 10).Value = 2020.2222
$ws.Cells.Item(27, 12).Value = 2020.2222
$ws.Cells.Item(27, 14).Value = -2234.2222
$ws.Cells.Item(40, 8).Value = 10071.5
$ws.Cells.Item(40, 9).Value = 11232.333
$ws.Cells.Item(40, 11).Value = 11232.333
$ws.Cells.Item(40, 13).Value = -11096.333
$ws.Cells.Item(93, 8).Value = 17544692
$ws.Cells.Item(93, 9).Value = 828.6429000000001
$ws.Cells.Item(93, 11).Value = 828.6429000000001
$ws.Cells.Item(93, 13).Value = 419.3570999999999
$ws.Cells.Item(122, 8).Value = 5072.769
$ws.Cells.Item(122, 9).Value = 4856.25
$ws.Cells.Item(122, 11).Value = 14568.75
$ws.Cells.Item(122, 13).Value = -12118.75
$ws.Cells.Item(126, 8).Value = 4074.9333
$ws.Cells.Item(126, 9).Value = 3306.875
$ws.Cells.Item(126, 10).Value = 4952.7144
$ws.Cells.Item(126, 11).Value = 9920.625
$ws.Cells.Item(126, 12).Value = 14858.1432
$ws.Cells.Item(126, 13).Value = -7450.625
$ws.Cells.Item(126, 14).Value = -19798.1432
$ws.Cells.Item(136, 8).Value = 3920.6428
$ws.Cells.Item(136, 9).Value = 3110.4736
$ws.Cells.Item(136, 10).Value = 5631
$ws.Cells.Item(136, 11).Value = 9331.4208
$ws.Cells.Item(136, 12).Value = 16893
$ws.Cells.Item(136, 13).Value = -6781.4208
$ws.Cells.Item(136, 14).Value = -21993

# --- Sheet: WVR (index 8) ---
$ws = $wb.Worksheets.Item(8)
$ws.Cells.Item(81, 8).Value = 1325.8334
$ws.Cells.Item(81, 9).Value = 591
$ws.Cells.Item(81, 11).Value = 1182
$ws.Cells.Item(81, 13).Value = -121
$ws.Cells.Item(84, 8).Value = 1325.8334
$ws.Cells.Item(84, 9).Value = 591
$ws.Cells.Item(84, 11).Value = 5910
$ws.Cells.Item(84, 13).Value = -606
$ws.Cells.Item(100, 8).Value = 701.2
$ws.Cells.Item(100, 9).Value = 526.5
$ws.Cells.Item(100, 10).Value = 1400
$ws.Cells.Item(100, 11).Value = 1053
$ws.Cells.Item(100, 12).Value = 2800
$ws.Cells.Item(100, 13).Value = -512
$ws.Cells.Item(100, 14).Value = -3882
$ws.Cells.Item(126, 8).Value = 16611.445
$ws.Cells.Item(126, 9).Value = 20450.5
$ws.Cells.Item(126, 11).Value = 61351.5
$ws.Cells.Item(126, 13).Value = -58881.5
$ws.Cells.Item(136, 8).Value = 10483447
$ws.Cells.Item(136, 9).Value = 16836278
$ws.Cells.Item(136, 10).Value = 1274.15
$ws.Cells.Item(136, 11).Value = 50508834
$ws.Cells.Item(136, 12).Value = 3822.45
$ws.Cells.Item(136, 13).Value = -50506284
$ws.Cells.Item(136, 14).Value = -8922.450000000001
